# Update the results table (tabela1) with the new data values.
# wdFindContinue = 1, wdReplaceAll = 2 (per the Execute signature used below)

$d = $word.ActiveDocument

# n: 426 -> 422
$d.Content.Find.Execute("426", $true, $false, $false, $false, $false, $true, 1, $false, "422", 2)

# Sexo = M (%): 159 (37.3) -> 159 (37.7)
$d.Content.Find.Execute("159 (37.3)", $true, $false, $false, $false, $false, $true, 1, $false, "159 (37.7)", 2)

# Idade (média (DP)): 63.43 (15.77) -> 63.43 (15.84)
$d.Content.Find.Execute("63.43 (15.77)", $true, $false, $false, $false, $false, $true, 1, $false, "63.43 (15.84)", 2)

# Peso (média (DP)): 76.29 (19.26) -> 76.28 (19.32)
$d.Content.Find.Execute("76.29 (19.26)", $true, $false, $false, $false, $false, $true, 1, $false, "76.28 (19.32)", 2)

# IMC (média (DP)): 28.71 (6.41) -> 28.68 (6.41)
$d.Content.Find.Execute("28.71 (6.41)", $true, $false, $false, $false, $false, $true, 1, $false, "28.68 (6.41)", 2)

# Tipo.Atendimento = INTERNAÇÃO (%): 250 (58.7) -> 250 (59.2)
$d.Content.Find.Execute("250 (58.7)", $true, $false, $false, $false, $false, $true, 1, $false, "250 (59.2)", 2)

# SPT (%): 54 (12.7) -> 54 (12.8)
$d.Content.Find.Execute("54 (12.7)", $true, $false, $false, $false, $false, $true, 1, $false, "54 (12.8)", 2)

# TVP.PREVIA (%): 62 (14.7) -> 61 (14.6)
$d.Content.Find.Execute("62 (14.7)", $true, $false, $false, $false, $false, $true, 1, $false, "61 (14.6)", 2)
